$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.003.83"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "2.586.63"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'521.74"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").Value = "'139.09"
$ws.Range("E6").Value = "  -2.20%  "
$ws.Range("D8").Value = "'0.564"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "2.597.88"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").Value = "'6.56"
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("D14").Value = "3.040.73"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "58.951.68"
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("D16").Value = "'20.44"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "2.599.66"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").Value = "'338.51"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").Value = "'4.29"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "'10.09"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("E22").Value = "  +3.07%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("D24").Value = "'65.95"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("D26").Value = "'0.403"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'7.01"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "0.0₃0724"
$ws.Range("E30").Value = "  -2.68%  "
$ws.Range("E31").Value = "  -4.69%  "
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "'18.68"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").Value = "'148.95"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").Value = "'3.98"
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("D37").Value = "'36.77"
$ws.Range("E37").Value = "  +2.31%  "
$ws.Range("E38").Value = "  +1.75%  "
$ws.Range("D39").Value = "'0.826"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").Value = "'0.816"
$ws.Range("E40").Value = "  -6.08%  "
$ws.Range("D41").Value = "'3.51"
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "'270.88"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("D44").Value = "'10.75"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("D45").Value = "'0.0953"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").Value = "'0.588"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("D47").Value = "'0.0517"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").Value = "'18.39"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").Value = "1.964.40"
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'4.51"
$ws.Range("E50").Value = "  -3.02%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0219"
$ws.Range("E51").Value = "  -0.28%  "
